# road map status updated
# Mark the "10-09-2025 / (0:36:47) Comments" row as DONE in the Status column
# of the roadmap table on slide 2.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tbl = $sh.Table

# Row 6, Column 3 is the "Status" cell (merged across the 2 rows covering
# 10-09-2025) that currently holds a single blank placeholder character.
$cell = $tbl.Cell(6, 3)
$cell.Shape.TextFrame.TextRange.Text = "DONE "
